$wb = $excel.ActiveWorkbook

# Sheet 1: 카카오 (Kakao)
$ws1 = $wb.Worksheets.Item("카카오")
$ws1.Range("B100").Value = 988877
$ws1.Range("A101").Value = 45960
$ws1.Range("B101").Value = 887895
$ws1.Range("A102").Value = 45961
$ws1.Range("B102").Value = 0
$ws1.Range("A102").NumberFormat = $ws1.Range("A101").NumberFormat

# Sheet 2: NAVER
$ws2 = $wb.Worksheets.Item("NAVER")
$ws2.Range("B100").Value = 1218200
$ws2.Range("A101").Value = 45960
$ws2.Range("B101").Value = 1155008
$ws2.Range("A102").Value = 45961
$ws2.Range("B102").Value = 0
$ws2.Range("A102").NumberFormat = $ws2.Range("A101").NumberFormat

# Sheet 3: 농심 (Nongshim)
$ws3 = $wb.Worksheets.Item("농심")
$ws3.Range("B100").Value = 138423
$ws3.Range("A101").Value = 45960
$ws3.Range("B101").Value = 135307
$ws3.Range("A102").Value = 45961
$ws3.Range("B102").Value = 0
$ws3.Range("A102").NumberFormat = $ws3.Range("A101").NumberFormat
